# Update the public EPEX Spot / Gaz / CO2 workbook with the latest daily data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column (CW) with header "22-sep" and the
# corresponding hourly price values for rows 2..25.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("CV1").Copy()
$wsPrix.Range("CW1").PasteSpecial(-4122)
$wsPrix.Range("CW1").Value = "22-sep"

$prixValues = @{
    2  = 15.57
    3  = 14.08
    4  = 21.44
    5  = 11.33
    6  = 4.14
    7  = 1.5
    8  = 10.66
    9  = 41.9
    10 = 62.03
    11 = 29.35
    12 = 2
    13 = 0
    14 = 0
    15 = -0.01
    16 = -0.01
    17 = -0.01
    18 = -0.01
    19 = -0.01
    20 = 1.54
    21 = 39.66
    22 = 17.28
    23 = 9.220000000000001
    24 = 6.7
    25 = 26.47
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 101).Value = $prixValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append rows for 2025-09-20 and 2025-09-21.
# Dates are stored as plain text (like the existing A column), so the
# number format is forced to Text before assignment and then cleared back
# to the sheet's default (unstyled) formatting.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A98").NumberFormat = "@"
$wsGaz.Range("A98").Value = "2025-09-20"
$wsGaz.Range("A98").ClearFormats()
$wsGaz.Range("B98").Value = 31.75

$wsGaz.Range("A99").NumberFormat = "@"
$wsGaz.Range("A99").Value = "2025-09-21"
$wsGaz.Range("A99").ClearFormats()
$wsGaz.Range("B99").Value = 31.75

# ---------------------------------------------------------------------------
# Sheet "CO2": append rows for 2025-09-20 and 2025-09-21.
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A98").NumberFormat = "@"
$wsCO2.Range("A98").Value = "2025-09-20"
$wsCO2.Range("A98").ClearFormats()
$wsCO2.Range("B98").Value = 76.63

$wsCO2.Range("A99").NumberFormat = "@"
$wsCO2.Range("A99").Value = "2025-09-21"
$wsCO2.Range("A99").ClearFormats()
$wsCO2.Range("B99").Value = 76.63
